$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Added ability to add new Todo ---
$ws.Range("A10").Value = 41843
$ws.Range("A10").NumberFormat = "mm-dd-yy"
$ws.Range("B10").Value = "Added ability to add new Todo"
$ws.Range("C10").Value = "Modal pops up with necessary fillings. More work should be done so that stuff like 'Email ATO' doesn't show up with 'Is Public' is not checked"
$ws.Rows.Item(10).RowHeight = 15.75

# --- Row 11: Set up login and cookies ---
$ws.Range("A11").Value = 41844
$ws.Range("A11").NumberFormat = "mm-dd-yy"
$ws.Range("B11").Value = "Set up login and cookies"
$ws.Range("C11").Value = "Users can now log in and it is saved as a cookie variable"
$ws.Rows.Item(11).RowHeight = 15.75

# Match the style (font/fill/border) already used for the Summary/Description columns
$ws.Range("B10:C11").Font.Name = "Arial"
$ws.Range("B10:C11").Font.Size = 10

# Widen column B to fit the new, longer "Summary" text (bestFit-style autosize)
$ws.Columns.Item(2).ColumnWidth = 25.6

# Move the active selection down past the newly added rows
$ws.Range("A12").Select() | Out-Null
